$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force these cells to remain text (values look numeric, Excel would
# otherwise auto-convert them to a Number type on assignment).
$textCells = @("D5", "D6", "D7", "D8", "D9", "D10", "D12", "D13", "D14", "D17", "D18", "D21", "D22", "D23", "D24", "D25", "D26", "D27", "D29", "D30", "D32", "D34", "D35", "D36", "D37", "D38", "D39", "D43", "D44", "D46", "D47", "D48", "D50", "D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated values row by row (Coin / Link / Price / Volume(1h)).
# Row 2
$ws.Range("D2").Value = "94.586.41"
$ws.Range("E2").Value = "  -3.27%  "
# Row 3
$ws.Range("D3").Value = "3.439.38"
$ws.Range("E3").Value = "  +2.46%  "
# Row 4
$ws.Range("E4").Value = "  -0.01%  "
# Row 5
$ws.Range("D5").Value = "237.66"
$ws.Range("E5").Value = "  -5.73%  "
# Row 6
$ws.Range("D6").Value = "643.38"
$ws.Range("E6").Value = "  -2.26%  "
# Row 7
$ws.Range("D7").Value = "1.44"
$ws.Range("E7").Value = "  +3.47%  "
# Row 8
$ws.Range("D8").Value = "0.406"
$ws.Range("E8").Value = "  -3.67%  "
# Row 9
$ws.Range("D9").Value = "0.999"
$ws.Range("E9").Value = "  +0.04%  "
# Row 10
$ws.Range("D10").Value = "0.974"
$ws.Range("E10").Value = "  -3.31%  "
# Row 11
$ws.Range("D11").Value = "3.438.68"
$ws.Range("E11").Value = "  +2.52%  "
# Row 12
$ws.Range("B12").Value = "TRON"
$ws.Range("C12").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D12").Value = "0.198"
$ws.Range("E12").Value = "  -4.72%  "
# Row 13
$ws.Range("B13").Value = "Avalanche"
$ws.Range("C13").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D13").Value = "41.82"
$ws.Range("E13").Value = "  +1.84%  "
# Row 14
$ws.Range("D14").Value = "6.22"
$ws.Range("E14").Value = "  +2.70%  "
# Row 15
$ws.Range("D15").Value = "94.271.10"
$ws.Range("E15").Value = "  -3.31%  "
# Row 16
$ws.Range("D16").Value = "4.075.21"
$ws.Range("E16").Value = "  +2.20%  "
# Row 17
$ws.Range("D17").Value = "0.0000252"
$ws.Range("E17").Value = "  -0.79%  "
# Row 18
$ws.Range("D18").Value = "8.36"
$ws.Range("E18").Value = "  -3.58%  "
# Row 19
$ws.Range("D19").Value = "3.431.19"
$ws.Range("E19").Value = "  +3.00%  "
# Row 20
$ws.Range("E20").Value = "  -0.83%  "
# Row 21
$ws.Range("D21").Value = "11.48"
$ws.Range("E21").Value = "  +5.73%  "
# Row 22
$ws.Range("D22").Value = "0.505"
$ws.Range("E22").Value = "  -1.38%  "
# Row 23
$ws.Range("D23").Value = "501.00"
$ws.Range("E23").Value = "  -1.61%  "
# Row 24
$ws.Range("D24").Value = "3.23"
$ws.Range("E24").Value = "  -3.46%  "
# Row 25
$ws.Range("D25").Value = "0.0000193"
$ws.Range("E25").Value = "  -3.49%  "
# Row 26
$ws.Range("D26").Value = "6.51"
$ws.Range("E26").Value = "  -5.80%  "
# Row 27
$ws.Range("D27").Value = "94.14"
$ws.Range("E27").Value = "  -2.22%  "
# Row 28
$ws.Range("D28").Value = "3.619.19"
$ws.Range("E28").Value = "  +2.44%  "
# Row 29
$ws.Range("D29").Value = "11.93"
$ws.Range("E29").Value = "  -2.07%  "
# Row 30
$ws.Range("D30").Value = "11.74"
$ws.Range("E30").Value = "  +3.92%  "
# Row 31
$ws.Range("E31").Value = "  -0.24%  "
# Row 32
$ws.Range("D32").Value = "2.76"
$ws.Range("E32").Value = "  +8.35%  "
# Row 33
$ws.Range("E33").Value = "  -1.49%  "
# Row 34
$ws.Range("D34").Value = "0.999"
$ws.Range("E34").Value = "  -0.64%  "
# Row 35
$ws.Range("D35").Value = "0.178"
$ws.Range("E35").Value = "  -3.46%  "
# Row 36
$ws.Range("D36").Value = "29.99"
$ws.Range("E36").Value = "  +5.45%  "
# Row 37
$ws.Range("D37").Value = "0.554"
$ws.Range("E37").Value = "  -0.33%  "
# Row 38
$ws.Range("D38").Value = "547.19"
$ws.Range("E38").Value = "  +4.54%  "
# Row 39
$ws.Range("D39").Value = "7.66"
$ws.Range("E39").Value = "  -4.26%  "
# Row 40
$ws.Range("E40").Value = "  -3.53%  "
# Row 41
$ws.Range("E41").Value = "  +0.04%  "
# Row 42
$ws.Range("E42").Value = "  -0.01%  "
# Row 43
$ws.Range("D43").Value = "0.908"
$ws.Range("E43").Value = "  +7.83%  "
# Row 44
$ws.Range("D44").Value = "24.06"
$ws.Range("E44").Value = "  -1.46%  "
# Row 45
$ws.Range("E45").Value = "  -1.47%  "
# Row 46
$ws.Range("B46").Value = "dogwifhat"
$ws.Range("C46").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D46").Value = "3.37"
$ws.Range("E46").Value = "  +6.73%  "
# Row 47
$ws.Range("B47").Value = "Filecoin"
$ws.Range("C47").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D47").Value = "5.58"
$ws.Range("E47").Value = "  +0.04%  "
# Row 48
$ws.Range("D48").Value = "0.0410"
$ws.Range("E48").Value = "  -3.10%  "
# Row 49
$ws.Range("E49").Value = "  -1.07%  "
# Row 50
$ws.Range("D50").Value = "54.41"
$ws.Range("E50").Value = "  -0.47%  "
# Row 51
$ws.Range("D51").Value = "2.17"
$ws.Range("E51").Value = "  -8.22%  "

# Restore default style on the forced-text cells (NumberFormat="@" alone
# would otherwise leave behind a new style index).
foreach ($addr in $textCells) {
    $ws.Range($addr).Style = "Normal"
}
